$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.207.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.25%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.271.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''299.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.06%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''95.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -4.29%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -2.15%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -0.05%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -2.79%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''33.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -4.94%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0787'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -0.65%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -6.64%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +1.78%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''16.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +2.09%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''6.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -0.54%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''2.625.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -1.58%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.260.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.75%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -2.10%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''42.151.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.25%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''11.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +1.00%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.0₃0889'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -1.88%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''5.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.40%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''66.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -2.36%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''235.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -0.01%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +0.21%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  +0.08%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -2.94%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''23.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -4.98%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''167.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +1.81%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -4.94%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''9.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.35%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''33.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -3.26%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = '''4.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +6.40%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -2.11%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''16.72'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.76%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''2.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -3.40%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.0688'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -1.78%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  -3.25%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.0987'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.37%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -1.78%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -4.39%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -6.80%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''1.959.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -0.43%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -0.94%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''17.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -5.06%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''9.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -6.08%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -3.88%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''2.496.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = '''52.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -6.11%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -3.79%  '
$ws.Range("E51").Style = "Normal"
